$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sCryo")

# Insert a new row at position 13 (pushes old row13 down to row14, etc.)
[void]$ws.Rows.Item(13).Insert()

# Fill new row 13 content
$ws.Cells.Item(13,1).Value = "cnsb"
$ws.Cells.Item(13,2).Value = "Snow sublimation"
$ws.Cells.Item(13,3).Value = "m / m^2"
$ws.Cells.Item(13,4).Value = "depth of water equivalent per unit area"
$ws.Cells.Item(13,5).Value = "Depth of snow (in water equivalent) that is lost due to sublimation. Needed for full water balance."
$ws.Rows.Item(13).RowHeight = 60

# Update Units column (C) for rows that referenced the old "m-1" shared string
# (index 160) to instead use the new "m / m^2" string, matching the new row's units.
$ws.Cells.Item(5,3).Value = "m / m^2"
$ws.Cells.Item(7,3).Value = "m / m^2"
$ws.Cells.Item(10,3).Value = "m / m^2"
$ws.Cells.Item(12,3).Value = "m / m^2"
$ws.Cells.Item(14,3).Value = "m / m^2"

[void]$ws.Range("D14").Select()

Write-Output "done"
